# Apply tracking updates to the tpivot phase-3 tracking sheet:
#  - Row 13 ("Bug: identical Value entries cause 'column ambiguously defined'...")
#    is now resolved. Its Approach text is replaced with the new resolution
#    note, status flips from "Not Started" to "Complete", completion dates are
#    recorded, and the row's "in progress" (green) highlight formatting is
#    cleared to match the other completed rows (e.g. row 15's formatting,
#    which already has the plain/no-highlight look with date columns).
#  - The active selection moves to E16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Pick up the "Complete"-style formatting (no highlight fill, date columns
# formatted, row height 38.25) from row 15 and apply it to row 13, without
# touching row 13's own values.
$ws.Range("A15:H15").Copy()
$ws.Range("A13:H13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Rows.Item(13).RowHeight = 38.25

# Update the row's content: new resolution/approach text, status, and dates.
$ws.Range("E13").Value2 = "Reject queries with duplicate filter signatures on server-side validation step. This just tosses the server request; there's no reason to return a query containing dupe values."
$ws.Range("F13").Value2 = "Complete"
$ws.Range("G13").Value2 = 43032
$ws.Range("H13").Value2 = 43032

# Update the active cell/selection on the sheet.
$ws.Range("E16").Select()
